# Sample Project / Main.xlsx - "Rules" sheet
# Cell B11 (row 11, the "R40" rule-label cell) is renamed from the text
# "R40" to the text "1" -- it must stay a literal text value (a shared
# string), not become a number, and it must keep its existing cell style.
#
# A plain  $ws.Range("B11").Value = "1"  assignment would be interpreted
# as the number 1 (since "1" parses as numeric), which also silently
# changes the cell's type. Forcing text via a leading apostrophe or via
# NumberFormat "@" does produce a text cell, but Excel then also flips on
# the cell's "quote prefix" bit, which mints a brand-new cell style -
# something the target workbook does not do (B11 keeps its original
# style index).
#
# The reliable way to land a literal, non-numeric-looking text value
# without touching the cell's style is to build the text on a scratch
# cell via a formula that evaluates to a string ("1"), copy it, and use
# Paste Special → Values onto B11: this overwrites only the value/type
# of B11 and leaves its formatting untouched. The scratch cell is then
# fully cleared so it leaves no trace in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scratch = $ws.Range("F3")
$scratch.Formula = "=""1"""
$scratch.Copy()
$ws.Range("B11").PasteSpecial(-4163)
$scratch.Clear()
